$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting all existing data
# (columns B..M) one column to the right (to C..N).
$ws.Range("B:B").Insert()

# Update the selection to match the post-edit state (row 3, which is
# the empty row sitting between the header row and the data rows).
$ws.Range("A3:XFD3").Select()
